$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking strings (e.g. "215.72")
# are preserved as text instead of being converted to numbers, matching
# the original inlineStr cell type.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '26.015.20'
$ws.Range('D3').Value = '1.639.96'
$ws.Range('E3').Value = '  -0.45%  '
$ws.Range('E4').Value = '  -0.59%  '
$ws.Range('D5').Value = '215.72'
$ws.Range('E5').Value = '  +0.16%  '
$ws.Range('D6').Value = '0.5158'
$ws.Range('E6').Value = '  +1.55%  '
$ws.Range('D7').Value = '1.002'
$ws.Range('E7').Value = '  -0.46%  '
$ws.Range('D8').Value = '0.2588'
$ws.Range('E8').Value = '  +0.56%  '
$ws.Range('E9').Value = '  -0.59%  '
$ws.Range('D10').Value = '19.88'
$ws.Range('E10').Value = '  +0.84%  '
$ws.Range('D11').Value = '0.07776'
$ws.Range('E11').Value = '  +0.05%  '
$ws.Range('D12').Value = '4.296'
$ws.Range('E12').Value = '  -0.33%  '
$ws.Range('D13').Value = '1.640.93'
$ws.Range('E13').Value = '  -0.75%  '
$ws.Range('D14').Value = '0.5479'
$ws.Range('E14').Value = '  +0.15%  '
$ws.Range('D15').Value = '0.0₅7789'
$ws.Range('E15').Value = '  -1.47%  '
$ws.Range('E16').Value = '  -0.85%  '
$ws.Range('D17').Value = '26.029.68'
$ws.Range('E17').Value = '  +0.12%  '
$ws.Range('E18').Value = '  -0.43%  '
$ws.Range('D19').Value = '199.22'
$ws.Range('E19').Value = '  +0.86%  '
$ws.Range('D20').Value = '4.463'
$ws.Range('E20').Value = '  +0.91%  '
$ws.Range('D21').Value = '9.989'
$ws.Range('E21').Value = '  -0.42%  '
$ws.Range('E22').Value = '  +0.70%  '
$ws.Range('D23').Value = '1.003'
$ws.Range('E23').Value = '  -0.58%  '
$ws.Range('E24').Value = '  +1.56%  '
$ws.Range('D25').Value = '142.27'
$ws.Range('E25').Value = '  +0.70%  '
$ws.Range('D26').Value = '0.1232'
$ws.Range('E26').Value = '  +7.55%  '
$ws.Range('E27').Value = '  -0.38%  '
$ws.Range('D28').Value = '15.62'
$ws.Range('E28').Value = '  -0.79%  '
$ws.Range('D29').Value = '1.244'
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('D30').Value = '0.04864'
$ws.Range('E30').Value = '  -3.30%  '
$ws.Range('D31').Value = '3.310'
$ws.Range('E31').Value = '  +1.02%  '
$ws.Range('E32').Value = '  +1.15%  '
$ws.Range('D33').Value = '1.543'
$ws.Range('E33').Value = '  +0.13%  '
$ws.Range('E34').Value = '  +0.43%  '
$ws.Range('D35').Value = '0.9217'
$ws.Range('E35').Value = '  +3.16%  '
$ws.Range('D36').Value = '0.5595'
$ws.Range('E36').Value = '  +1.01%  '
$ws.Range('D37').Value = '2.569'
$ws.Range('E37').Value = '  -0.99%  '
$ws.Range('D38').Value = '1.114.88'
$ws.Range('E38').Value = '  -1.69%  '
$ws.Range('D39').Value = '0.01575'
$ws.Range('E39').Value = '  +0.71%  '
$ws.Range('E40').Value = '  -0.61%  '
$ws.Range('D41').Value = '2.534'
$ws.Range('E41').Value = '  -0.95%  '
$ws.Range('D42').Value = '5.575'
$ws.Range('E42').Value = '  -1.68%  '
$ws.Range('D43').Value = '0.8093'
$ws.Range('E43').Value = '  -0.74%  '
$ws.Range('D44').Value = '99.70'
$ws.Range('E44').Value = '  -0.10%  '
$ws.Range('E45').Value = '  -0.16%  '
$ws.Range('D46').Value = '1.779.94'
$ws.Range('E46').Value = '  -0.31%  '
$ws.Range('E47').Value = '  -0.03%  '
$ws.Range('D48').Value = '55.37'
$ws.Range('E48').Value = '  +0.31%  '
$ws.Range('E49').Value = '  +0.11%  '
$ws.Range('D50').Value = '0.05213'
$ws.Range('E50').Value = '  +2.34%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '7.530'
$ws.Range('E51').Value = '  +1.72%  '

# Restore column D styling back to the workbook default ("Normal") so the
# cells do not retain an extraneous explicit style index.
$ws.Range("D2:D51").Style = "Normal"

